$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1201
$ws.Range("I28").Value = 944.7143
$ws.Range("J28").Value = 2995
$ws.Range("K28").Value = 944.7143
$ws.Range("L28").Value = 2995
$ws.Range("M28").Value = -459.7143
$ws.Range("N28").Value = -3965

# Row 98
$ws.Range("H98").Value = 870.5
$ws.Range("I98").Value = 305.75
$ws.Range("K98").Value = 305.75
$ws.Range("M98").Value = 1192.25

# Row 100
$ws.Range("I100").Value = 2981.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2981.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2440.5
$ws.Range("N100").Value = -4082

# Row 103
$ws.Range("H103").Value = 999.75
$ws.Range("I103").Value = 999
$ws.Range("K103").Value = 2997
$ws.Range("M103").Value = -2411

# Row 107
$ws.Range("H107").Value = 281.1111
$ws.Range("I107").Value = 253.75
$ws.Range("K107").Value = 253.75
$ws.Range("M107").Value = 1666.25

# Row 116
$ws.Range("H116").Value = 3414.6155
$ws.Range("I116").Value = 2346.25
$ws.Range("K116").Value = 2346.25
$ws.Range("M116").Value = 1095.75

# Row 122
$ws.Range("H122").Value = 870.5
$ws.Range("I122").Value = 305.75
$ws.Range("K122").Value = 917.25
$ws.Range("M122").Value = 1532.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 844
$ws.Range("I2").Value = 818.8461
$ws.Range("J2").Value = 884.875
$ws.Range("K2").Value = 818.8461
$ws.Range("L2").Value = 884.875
$ws.Range("M2").Value = -705.8461
$ws.Range("N2").Value = -1110.875

# Row 102
$ws.Range("H102").Value = 1654.5
$ws.Range("I102").Value = 1645.4
$ws.Range("K102").Value = 1645.4
$ws.Range("M102").Value = -23.40000000000009

# Row 116
$ws.Range("H116").Value = 844
$ws.Range("I116").Value = 818.8461
$ws.Range("J116").Value = 884.875
$ws.Range("K116").Value = 818.8461
$ws.Range("L116").Value = 884.875
$ws.Range("M116").Value = 1475.1539
$ws.Range("N116").Value = -5472.875

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 844
$ws.Range("I3").Value = 818.8461
$ws.Range("J3").Value = 884.875
$ws.Range("K3").Value = 818.8461
$ws.Range("L3").Value = 884.875
$ws.Range("M3").Value = -704.8461
$ws.Range("N3").Value = -1112.875

$ws = $wb.Worksheets.Item("CRP")
# Row 92
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

# Row 122
$ws.Range("H122").Value = 2291.6428
$ws.Range("I122").Value = 1738.2
$ws.Range("K122").Value = 5214.6
$ws.Range("M122").Value = -2764.6

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 1797.6666
$ws.Range("I33").Value = 197.25
$ws.Range("K33").Value = 1183.5
$ws.Range("M33").Value = -900.5

# Row 92
$ws.Range("H92").Value = 570
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 683.3333
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 2049.9999
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -4545.9999

# Row 113
$ws.Range("H113").Value = 1044.125
$ws.Range("I113").Value = 724.5
$ws.Range("J113").Value = 1089.7858
$ws.Range("K113").Value = 2173.5
$ws.Range("L113").Value = 3269.3574
$ws.Range("M113").Value = -3.5
$ws.Range("N113").Value = -7609.357400000001

# Row 131
$ws.Range("H131").Value = 1429.6842
$ws.Range("I131").Value = 959.5714
$ws.Range("K131").Value = 2878.7142
$ws.Range("M131").Value = 2161.2858

# Row 132
$ws.Range("H132").Value = 4264.4443
$ws.Range("I132").Value = 3745.75
$ws.Range("K132").Value = 33711.75
$ws.Range("M132").Value = -31181.75

# Row 137
$ws.Range("H137").Value = 3149.5
$ws.Range("I137").Value = 3300
$ws.Range("J137").Value = 2999
$ws.Range("K137").Value = 9900
$ws.Range("L137").Value = 8997
$ws.Range("M137").Value = -4800
$ws.Range("N137").Value = -19197

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 22962
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 27702.5
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 27702.5
$ws.Range("M46").Value = -3844
$ws.Range("N46").Value = -28014.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 924.25
$ws.Range("I46").Value = 899.3333
$ws.Range("K46").Value = 899.3333
$ws.Range("M46").Value = -711.3333

# Row 93
$ws.Range("H93").Value = 2842.5
$ws.Range("I93").Value = 2827.6667
$ws.Range("J93").Value = 2887
$ws.Range("K93").Value = 2827.6667
$ws.Range("L93").Value = 2887
$ws.Range("M93").Value = -1579.6667
$ws.Range("N93").Value = -5383

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 13364.272
$ws.Range("I81").Value = 7833
$ws.Range("J81").Value = 20001.8
$ws.Range("K81").Value = 15666
$ws.Range("L81").Value = 40003.6
$ws.Range("M81").Value = -14605
$ws.Range("N81").Value = -42125.6

# Row 84
$ws.Range("H84").Value = 13364.272
$ws.Range("I84").Value = 7833
$ws.Range("J84").Value = 20001.8
$ws.Range("K84").Value = 78330
$ws.Range("L84").Value = 200018
$ws.Range("M84").Value = -73026
$ws.Range("N84").Value = -210626

# Row 100
$ws.Range("H100").Value = 1175.7142
$ws.Range("I100").Value = 1205
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 2410
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1869
$ws.Range("N100").Value = -3082

# Row 122
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 3173.9
$ws.Range("I132").Value = 2552.1177
$ws.Range("J132").Value = 3987
$ws.Range("K132").Value = 7656.353099999999
$ws.Range("L132").Value = 11961
$ws.Range("M132").Value = -5126.353099999999
$ws.Range("N132").Value = -17021

# Row 136
$ws.Range("H136").Value = 2292.8
$ws.Range("I136").Value = 1798.625
$ws.Range("J136").Value = 2857.5715
$ws.Range("K136").Value = 5395.875
$ws.Range("L136").Value = 8572.7145
$ws.Range("M136").Value = -2845.875
$ws.Range("N136").Value = -13672.7145
